# Update countries & provincias Spain
# Applies the COVID data refresh: updated figures for several countries,
# re-ranking of a few countries (Uzbekistan, Costa Rica, Maldivas move up
# in the ranking, pushing the countries below them down one row), and an
# updated "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4,2).Value = 4356238
$ws.Cells.Item(4,3).Value = 40529
$ws.Cells.Item(4,4).Value = 2075495
$ws.Cells.Item(4,5).Value = 2131107
$ws.Cells.Item(4,7).Value = 238
$ws.Cells.Item(4,8).Value = 149636

# --- Row 6: India ---
$ws.Cells.Item(6,2).Value = 1436019
$ws.Cells.Item(6,3).Value = 50525
$ws.Cells.Item(6,4).Value = 918735

# --- Row 8: Sudafrica ---
$ws.Cells.Item(8,2).Value = 445433
$ws.Cells.Item(8,3).Value = 11233
$ws.Cells.Item(8,4).Value = 265077
$ws.Cells.Item(8,5).Value = 173587
$ws.Cells.Item(8,7).Value = 114
$ws.Cells.Item(8,8).Value = 6769

# --- Rows 64-66: Uzbekistan overtakes Marruecos and Austria ---
# Row 64 becomes Uzbekistan with refreshed figures
$ws.Cells.Item(64,1).Value = "Uzbekistan"
$ws.Cells.Item(64,2).Value = 20531
$ws.Cells.Item(64,3).Value = 579
$ws.Cells.Item(64,4).Value = 11105
$ws.Cells.Item(64,5).Value = 9310
$ws.Cells.Item(64,7).Value = 5
$ws.Cells.Item(64,8).Value = 116

# Row 65 becomes Austria (previous row 64 data)
$ws.Cells.Item(65,1).Value = "Austria"
$ws.Cells.Item(65,2).Value = 20472
$ws.Cells.Item(65,3).Value = 134
$ws.Cells.Item(65,4).Value = 18209
$ws.Cells.Item(65,5).Value = 1551
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 712

# Row 66 becomes Marruecos (previous row 65 data)
$ws.Cells.Item(66,1).Value = "Marruecos"
$ws.Cells.Item(66,2).Value = 20278
$ws.Cells.Item(66,3).Value = 633
$ws.Cells.Item(66,4).Value = 16438
$ws.Cells.Item(66,5).Value = 3527
$ws.Cells.Item(66,7).Value = 8
$ws.Cells.Item(66,8).Value = 313

# --- Rows 72-74: Costa Rica overtakes El Salvador and Venezuela ---
# Row 72 becomes Costa Rica with refreshed figures
$ws.Cells.Item(72,1).Value = "Costa Rica"
$ws.Cells.Item(72,2).Value = 15229
$ws.Cells.Item(72,3).Value = 629
$ws.Cells.Item(72,4).Value = 3736
$ws.Cells.Item(72,5).Value = 11389
$ws.Cells.Item(72,7).Value = 6
$ws.Cells.Item(72,8).Value = 104

# Row 73 becomes Venezuela (previous row 72 data)
$ws.Cells.Item(73,1).Value = "Venezuela"
$ws.Cells.Item(73,2).Value = 14929
$ws.Cells.Item(73,3).Value = 0
$ws.Cells.Item(73,4).Value = 8795
$ws.Cells.Item(73,5).Value = 5996
$ws.Cells.Item(73,7).Value = 0
$ws.Cells.Item(73,8).Value = 138

# Row 74 becomes El Salvador (previous row 73 data)
$ws.Cells.Item(74,1).Value = "El Salvador"
$ws.Cells.Item(74,2).Value = 14630
$ws.Cells.Item(74,3).Value = 409
$ws.Cells.Item(74,4).Value = 7648
$ws.Cells.Item(74,5).Value = 6582
$ws.Cells.Item(74,7).Value = 10
$ws.Cells.Item(74,8).Value = 400

# --- Rows 108-109: Maldivas overtakes Tailandia ---
# Row 108 becomes Maldivas with refreshed figures
$ws.Cells.Item(108,1).Value = "Maldivas"
$ws.Cells.Item(108,2).Value = 3302
$ws.Cells.Item(108,3).Value = 50
$ws.Cells.Item(108,4).Value = 2528
$ws.Cells.Item(108,5).Value = 759
$ws.Cells.Item(108,8).Value = 15

# Row 109 becomes Tailandia (previous row 108 data)
$ws.Cells.Item(109,1).Value = "Tailandia"
$ws.Cells.Item(109,2).Value = 3291
$ws.Cells.Item(109,3).Value = 9
$ws.Cells.Item(109,4).Value = 3109
$ws.Cells.Item(109,5).Value = 124
$ws.Cells.Item(109,8).Value = 58

# --- Row 148: Angola ---
$ws.Cells.Item(148,2).Value = 932
$ws.Cells.Item(148,3).Value = 16
$ws.Cells.Item(148,5).Value = 650
$ws.Cells.Item(148,7).Value = 1
$ws.Cells.Item(148,8).Value = 40

# --- Updated timestamp label ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 22:13"
